$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 'Rappahannock County High School'
$ws.Range("C9").Value = 'education'
$ws.Range("D9").Value = '12576 Lee Hwy., Washington, VA 22747'
$ws.Range("G9").Value = -81.7352471
$ws.Range("H9").Value = 36.7832503

$ws.Range("B10").Value = 'Rappahannock Elementary School'
$ws.Range("C10").Value = 'education'
$ws.Range("D10").Value = '34 School House Rd, Washington, VA 22747'
$ws.Range("G10").Value = -78.1758117675781
$ws.Range("H10").Value = 38.6886138916015

$ws.Range("B11").Value = 'BABY BEAR DAYCARE'
$ws.Range("C11").Value = 'education'
$ws.Range("D11").Value = '30 Main St., Sperryville (540) ­987­-9644'
$ws.Range("G11").Value = -78.2302489
$ws.Range("H11").Value = 38.65898

$ws.Range("B12").Value = 'HEADWATERS'
$ws.Range("C12").Value = 'education'
$ws.Range("D12").Value = '567 Mount Salem Avenue, Washington VA 22747'
$ws.Range("G12").Value = -78.1580047607421
$ws.Range("H12").Value = 38.7095642089843

$ws.Range("B13").Value = 'HEARTHSTONE SCHOOL'
$ws.Range("C13").Value = 'education'
$ws.Range("D13").Value = '11576 Lee Highway P.O. Box 247 Sperryville, Virginia 22740 (540) ­987-­9212'
$ws.Range("G13").Value = -78.2182748
$ws.Range("H13").Value = 38.6626933

$ws.Range("B14").Value = 'MOUNTAIN LAUREL MONTESSORI FARM SCHOOL'
$ws.Range("C14").Value = 'education'
$ws.Range("D14").Value = '23 Sunny Slope Ln., Flint Hill (540) ­675-­1011'
$ws.Range("G14").Value = -78.1042098999023
$ws.Range("H14").Value = 38.772533416748

$ws.Range("B15").Value = 'RAPPAHANNOCK CENTER FOR EDUCATION'
$ws.Range("C15").Value = 'education'
$ws.Range("D15").Value = '12018 Lee Highway, Sperryville VA 22740'
$ws.Range("G15").Value = -78.2227401733398
$ws.Range("H15").Value = 38.6597290039062

$ws.Range("B16").Value = 'RAPPAHANNOCK COUNTY PUBLIC SCHOOLS'
$ws.Range("C16").Value = 'education'
$ws.Range("D16").Value = '6 Schoolhouse Rd., Washington, VA 22747 (540) ­227­-0023'
$ws.Range("G16").Value = -81.63742
$ws.Range("H16").Value = 36.667606

$ws.Range("B17").Value = 'THE CHILD CARE & LEARNING CENTER'
$ws.Range("C17").Value = 'education'
$ws.Range("D17").Value = '12763 Lee Hwy Washington, VA 22747'
$ws.Range("G17").Value = -81.7352471
$ws.Range("H17").Value = 36.7832503

$ws.Range("B18").Value = 'RUTH''S ROOTS RESEARCH'
$ws.Range("C18").Value = 'education'
$ws.Range("D18").Value = '59 Bunker Hill Ln, Castleton, VA 22716 (540) 229-2225'
$ws.Range("G18").Value = -78.125987
$ws.Range("H18").Value = 38.643961

$ws.Range("B19").Value = 'WAKEFIELD COUNTRY DAY SCHOOL'
$ws.Range("C19").Value = 'education'
$ws.Range("D19").Value = '1059 Zachary Taylor Hwy Huntly, Virginia 22640'
$ws.Range("G19").Value = -78.1045379638671
$ws.Range("H19").Value = 38.7916221618652

# Column widths (bestFit applied by the original author via AutoFit in real Excel)
$ws.Columns.Item(2).ColumnWidth = 32.45
$ws.Columns.Item(4).ColumnWidth = 62.3
$ws.Columns.Item(6).ColumnWidth = 4.8

# Selection + print setup
[void]$ws.Range("C23").Select()
$ws.PageSetup.Orientation = 1

